$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (Neutro / Hc / C5ar1 -> M1) ---
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 0.1527673333333333
$ws.Cells.Item(2, 8).Value = 0.458302
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.169425
$ws.Cells.Item(2, 14).Value = 0.508275
$ws.Cells.Item(2, 15).Value = 0.0005062042177611658
$ws.Cells.Item(2, 16).Value = 0.0005062042177611658
$ws.Cells.Item(2, 17).Value = 0.02588260545000001
$ws.Cells.Item(2, 18).Value = 0.23294344905
$ws.Cells.Item(2, 19).Value = 0.0005062042177611658
$ws.Cells.Item(2, 20).Value = 0.0005062042177611658

# --- Update existing row 3 (Neutro / Hc / C5ar1 -> M2) ---
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 0.1527673333333333
$ws.Cells.Item(3, 8).Value = 0.458302
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 69.07203766666666
$ws.Cells.Item(3, 14).Value = 207.216113
$ws.Cells.Item(3, 15).Value = 0.2063718860630059
$ws.Cells.Item(3, 16).Value = 0.2063718860630059
$ws.Cells.Item(3, 17).Value = 10.55195100223622
$ws.Cells.Item(3, 18).Value = 94.967559020126
$ws.Cells.Item(3, 19).Value = 0.2063718860630059
$ws.Cells.Item(3, 20).Value = 0.2063718860630059

# --- Insert a new row at position 4 for the FAPs target cluster ---
$ws.Rows.Item(4).Insert()

$ws.Cells.Item(4, 1).Value = "Neutro"
$ws.Cells.Item(4, 2).Value = "Hc"
$ws.Cells.Item(4, 3).Value = "C5ar1"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.1527673333333333
$ws.Cells.Item(4, 8).Value = 0.458302
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 31.266034
$ws.Cells.Item(4, 14).Value = 93.798102
$ws.Cells.Item(4, 15).Value = 0.0934159556350244
$ws.Cells.Item(4, 16).Value = 0.09341595563502442
$ws.Cells.Item(4, 17).Value = 4.776428638089333
$ws.Cells.Item(4, 18).Value = 42.987857742804
$ws.Cells.Item(4, 19).Value = 0.0934159556350244
$ws.Cells.Item(4, 20).Value = 0.09341595563502442

# --- Update the row that was previously row 4 (now shifted to row 5) with new values ---
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 0.1527673333333333
$ws.Cells.Item(5, 8).Value = 0.458302
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 234.189438
$ws.Cells.Item(5, 14).Value = 702.5683140000001
$ws.Cells.Item(5, 15).Value = 0.6997059540842084
$ws.Cells.Item(5, 16).Value = 0.6997059540842085
$ws.Cells.Item(5, 17).Value = 35.776495938092
$ws.Cells.Item(5, 18).Value = 321.988463442828
$ws.Cells.Item(5, 19).Value = 0.6997059540842084
$ws.Cells.Item(5, 20).Value = 0.6997059540842085
